$wb = $excel.ActiveWorkbook

# --- "config" sheet (sheet1.xml) ---
$ws = $wb.Worksheets.Item("config")

# B5: iterations value change 1814400 -> 360000
$ws.Range("B5").Value = 360000

# New row 11: Log File / log_entries.txt / Gen1MaxP / Gen*FuelCons / GenConfig[3-6]
$ws.Range("A11").Value = "Log File"
$ws.Range("B11").Value = "log_entries.txt"
$ws.Range("C11").Value = "Gen1MaxP"
$ws.Range("D11").Value = "Gen*FuelCons"
$ws.Range("E11").Value = "GenConfig[3-6]"

# Update the data validation list attached to A11 to include "Template"
$ws.Range("A11").Validation.Modify(3, 1, 1, '",,,,,Community Name,Template,input,output"')

# --- "GenStats" sheet (sheet2.xml): move selection ---
$ws2 = $wb.Worksheets.Item("GenStats")
[void]$ws2.Range("K1").Select()

# --- "GenConfigurations" sheet (sheet4.xml): add a selection ---
$ws4 = $wb.Worksheets.Item("GenConfigurations")
[void]$ws4.Range("D1").Select()

# Re-activate "config" sheet and move its selection to C23 (keeps it the tab shown when opened)
[void]$ws.Activate()
[void]$ws.Range("C23").Select()
